# Refresh the crypto price / 1h-volume snapshot, and fix the two rows
# (Cardano/OKB and EthereumClassic/ImmutableX) whose rank order swapped.
#
# Column D holds prices as plain text (e.g. "27.442.92", "336.54") in the
# source data. Several of the new values look like ordinary decimal
# numbers (e.g. "336.54"), and Excel auto-converts a bare numeric-looking
# .Value assignment into a Double, which would lose the original text
# formatting (e.g. trailing zeros, "0.00001106" turning into scientific
# notation). To keep those cells text, exactly like the source file, we
# write them with a leading apostrophe (Excels force-text marker) and
# then reset the cell style to Normal so no stray quote-prefix formatting
# is left behind.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.431.10'
$ws.Range("E2").Value = '  +3.57%  '

$ws.Range("D3").Value = '1.798.42'
$ws.Range("E3").Value = '  +4.48%  '

$ws.Range("D4").Value = '''1.006'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.36%  '

$ws.Range("D5").Value = '''336.54'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.35%  '

$ws.Range("D6").Value = '''1.003'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.39%  '

$ws.Range("D7").Value = '''0.3792'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +2.44%  '

$ws.Range("B8").Value = 'OKB'
$ws.Range("C8").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D8").Value = '''49.19'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +2.11%  '

$ws.Range("B9").Value = 'Cardano'
$ws.Range("C9").Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range("D9").Value = '''0.3462'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +3.19%  '

$ws.Range("D10").Value = '''1.211'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +2.32%  '

$ws.Range("D11").Value = '''0.07565'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +2.34%  '

$ws.Range("E12").Value = '  +0.60%  '

$ws.Range("D13").Value = '''21.98'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +9.47%  '

$ws.Range("D14").Value = '''6.530'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +2.19%  '

$ws.Range("D15").Value = '1.794.93'
$ws.Range("E15").Value = '  +4.34%  '

$ws.Range("D16").Value = '''7.067'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.16%  '

$ws.Range("D17").Value = '''0.00001106'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +3.60%  '

$ws.Range("D18").Value = '''0.06677'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.85%  '

$ws.Range("D19").Value = '''84.92'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +3.51%  '

$ws.Range("D20").Value = '''1.003'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.30%  '

$ws.Range("D21").Value = '''17.48'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +5.83%  '

$ws.Range("D22").Value = '''6.486'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +5.40%  '

$ws.Range("D23").Value = '27.417.90'
$ws.Range("E23").Value = '  +3.58%  '

$ws.Range("D24").Value = '''12.55'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -1.84%  '

$ws.Range("D25").Value = '''2.454'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.67%  '

$ws.Range("D26").Value = '''2.577'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +7.68%  '

$ws.Range("B27").Value = 'ImmutableX'
$ws.Range("C27").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D27").Value = '''1.503'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +7.02%  '

$ws.Range("B28").Value = 'EthereumClassic'
$ws.Range("C28").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D28").Value = '''21.52'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +11.06%  '

$ws.Range("D29").Value = '''150.38'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.08%  '

$ws.Range("D30").Value = '2.000.71'
$ws.Range("E30").Value = '  +4.76%  '

$ws.Range("D31").Value = '''133.89'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +2.42%  '

$ws.Range("D32").Value = '''4.099'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.50%  '

$ws.Range("D33").Value = '''6.134'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +3.03%  '

$ws.Range("D34").Value = '''0.08701'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +1.14%  '

$ws.Range("D35").Value = '''13.29'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +5.04%  '

$ws.Range("D36").Value = '''1.679'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.73%  '

$ws.Range("D37").Value = '''5.490'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +3.02%  '

$ws.Range("D38").Value = '''0.6890'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +11.92%  '

$ws.Range("D39").Value = '''0.2212'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +2.84%  '

$ws.Range("D40").Value = '''0.02363'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +1.78%  '

$ws.Range("D41").Value = '''0.06359'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +2.84%  '

$ws.Range("D42").Value = '''8.828'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +4.61%  '

$ws.Range("D43").Value = '''1.274'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +4.44%  '

$ws.Range("D44").Value = '''14.42'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +2.20%  '

$ws.Range("D45").Value = '''0.6467'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +8.71%  '

$ws.Range("D46").Value = '''1.002'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.37%  '

$ws.Range("D47").Value = '''3.851'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.38%  '

$ws.Range("D48").Value = '''2.135'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +5.16%  '

$ws.Range("E49").Value = '  +2.17%  '

$ws.Range("D50").Value = '''0.07226'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.97%  '

$ws.Range("D51").Value = '''79.84'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +3.77%  '
